$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.28919917345047
$ws.Range("B1").Value = 0.213948518037796
$ws.Range("C1").Value = 1.686006307601929
$ws.Range("D1").Value = 3.748922348022461
$ws.Range("E1").Value = 3.078108072280884
